$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Set every IP Address cell in column C (rows 4-16) to the same value
for ($r = 4; $r -le 16; $r++) {
    $ws.Cells.Item($r, 3).Value = "54.90.200.200"
}

# Update the active selection to match the saved view state
$ws.Range("G11").Select()
